$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-13 with the new (shifted) year values
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 6933.38
$ws.Range("C2").Value = 31259
$ws.Range("D2").Value = 151.372972
$ws.Range("E2").Value = 21.63
$ws.Range("F2").Value = 20768.73
$ws.Range("G2").Value = 164.328062

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 6933.38
$ws.Range("C3").Value = 31259
$ws.Range("D3").Value = 151.372972
$ws.Range("E3").Value = 21.63
$ws.Range("F3").Value = 20768.73
$ws.Range("G3").Value = 164.328062

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 6933.38
$ws.Range("C4").Value = 31259
$ws.Range("D4").Value = 151.372972
$ws.Range("E4").Value = 21.63
$ws.Range("F4").Value = 20768.73
$ws.Range("G4").Value = 164.328062

$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 6933.38
$ws.Range("C5").Value = 31259
$ws.Range("D5").Value = 151.372972
$ws.Range("E5").Value = 21.63
$ws.Range("F5").Value = 20768.73
$ws.Range("G5").Value = 164.328062

$ws.Range("A6").Value = "2014年"
$ws.Range("B6").Value = 8003.1
$ws.Range("C6").Value = 32368.55
$ws.Range("D6").Value = 175.602299
$ws.Range("E6").Value = 22.96
$ws.Range("F6").Value = 22044.62
$ws.Range("G6").Value = 190.07132

$ws.Range("A7").Value = "2015年"
$ws.Range("B7").Value = 8003.1
$ws.Range("C7").Value = 32368.55
$ws.Range("D7").Value = 175.602299
$ws.Range("E7").Value = 22.96
$ws.Range("F7").Value = 22044.62
$ws.Range("G7").Value = 190.07132

$ws.Range("A8").Value = "2016年"
$ws.Range("B8").Value = 8003.1
$ws.Range("C8").Value = 32368.55
$ws.Range("D8").Value = 175.602299
$ws.Range("E8").Value = 22.96
$ws.Range("F8").Value = 22044.62
$ws.Range("G8").Value = 190.07132

$ws.Range("A9").Value = "2017年"
$ws.Range("B9").Value = 8003.1
$ws.Range("C9").Value = 32368.55
$ws.Range("D9").Value = 175.602299
$ws.Range("E9").Value = 22.96
$ws.Range("F9").Value = 22044.62
$ws.Range("G9").Value = 190.07132

$ws.Range("A10").Value = "2018年"
$ws.Range("B10").Value = 8003.1
$ws.Range("C10").Value = 32368.55
$ws.Range("D10").Value = 175.602299
$ws.Range("E10").Value = 22.96
$ws.Range("F10").Value = 22044.62
$ws.Range("G10").Value = 190.07132

$ws.Range("A11").Value = "2019年"
$ws.Range("B11").Value = 8003.1
$ws.Range("C11").Value = 32368.55
$ws.Range("D11").Value = 175.602299
$ws.Range("E11").Value = 22.96
$ws.Range("F11").Value = 22044.62
$ws.Range("G11").Value = 190.07132

$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = 8003.1
$ws.Range("C12").Value = 32368.55
$ws.Range("D12").Value = 175.602299
$ws.Range("E12").Value = 22.96
$ws.Range("F12").Value = 22044.62
$ws.Range("G12").Value = 190.07132

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 8003.1
$ws.Range("C13").Value = 32368.55
$ws.Range("D13").Value = 175.602299
$ws.Range("E13").Value = 22.96
$ws.Range("F13").Value = 22044.62
$ws.Range("G13").Value = 190.07132

# Remove the now-obsolete trailing rows (14-18)
$ws.Rows("14:18").Delete()
